# "Generate Report for Archive"
#
# The handoff/handback status for both locales has moved on from
# "Ready for handoff" to "In Translation" everywhere it is reported:
#   - Overview sheet: per-locale status columns (zh-cn / de-de), rows 2-3
#   - zh-cn sheet: Status column, rows 2-3
#   - de-de sheet: Status column, rows 2-3
# Afterwards the now-shorter status text lets the Status column(s) narrow
# to fit the content again.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# Update every occurrence of the old status text on each sheet.
$wsOverview.Cells.Replace($oldStatus, $newStatus) | Out-Null
$wsZhCn.Cells.Replace($oldStatus, $newStatus)     | Out-Null
$wsDeDe.Cells.Replace($oldStatus, $newStatus)     | Out-Null

# Re-fit the status columns now that the text is shorter: "Overview"
# reports status in columns E (zh-cn) and F (de-de); the locale sheets
# report it in column C.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth     = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth     = 12.5
